$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "305.70"
    "E2" = "0.19%"
    "E3" = "-0.80%"
    "D4" = "5.060"
    "E4" = "0.55%"
    "D5" = "0.07865"
    "E5" = "0.23%"
    "D6" = "2.211"
    "E6" = "2.58%"
    "D7" = "7.978"
    "E7" = "-0.75%"
    "D8" = "0.9281"
    "E8" = "0.61%"
    "D9" = "0.09773"
    "E9" = "-1.60%"
    "D10" = "0.1868"
    "E10" = "-0.07%"
    "D11" = "0.08966"
    "E11" = "2.53%"
    "D12" = "0.03783"
    "E12" = "4.89%"
    "D13" = "0.09901"
    "E13" = "-0.20%"
    "D14" = "0.001445"
    "E14" = "-1.89%"
    "D15" = "0.005677"
    "E15" = "0.15%"
    "E16" = "0.24%"
    "D17" = "4.154"
    "E17" = "2.52%"
    "E18" = "14.01%"
    "E19" = "-0.79%"
    "D20" = "0.1322"
    "E20" = "-1.77%"
    "D21" = "5.148"
    "E21" = "4.56%"
    "E22" = "2.57%"
    "D23" = "0.04590"
    "E23" = "-0.27%"
    "D24" = "0.001235"
    "E24" = "0.20%"
    "D25" = "0.004782"
    "E25" = "-7.70%"
    "E26" = "-6.76%"
    "E27" = "74.31%"
    "D39" = "0.01929"
    "E39" = "6.45%"
    "D40" = "0.04968"
    "E40" = "4.65%"
    "D41" = "0.007804"
    "E41" = "-0.79%"
    "D42" = "0.1391"
    "E42" = "-0.93%"
    "D43" = "0.007833"
    "E43" = "3.11%"
    "D44" = "0.002211"
    "E44" = "-0.50%"
    "E45" = "7.82%"
    "D46" = "0.00006293"
    "E46" = "-1.09%"
    "E47" = "0.29%"
    "E48" = "0.16%"
    "D49" = "51.76"
    "E49" = "49.95%"
    "D50" = "0.001907"
    "E50" = "-29.21%"
    "D51" = "0.00002107"
    "E51" = "0.29%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
